# Add 2022-Q3 data
# 1) Insert a brand-new worksheet "2022-Q3" right before "2022-Q2"
#    (so the tab order becomes 总计, 2022-Q3, 2022-Q2, 2022-Q1, ...)
# 2) Populate it with the fund-holding detail rows for 2022-Q3
# 3) Insert a matching summary row at the top of the "总计" sheet's data

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q3" worksheet, positioned before "2022-Q2"
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2)
$newSheet.Name = "2022-Q3"

# Re-fetch references by name after the structural change above so we
# never hold a stale sheet object.
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------
# Step 2: copy header-row + data-row formatting from the 2022-Q2 sheet
#         (same layout for every quarterly detail sheet) then fill values
# ---------------------------------------------------------------------
$q2.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$q2.Range("A2").Copy()
$q3.Range("A2:A3").PasteSpecial(-4122)

# NOTE: this runtime's `.Value` getter/setter is unreliable (it can hand
# back a stringified property descriptor instead of the real value), so
# every read/write below goes through `.Value2`, which behaves correctly.
$q3.Range("B1").Value2 = "基金代码"
$q3.Range("C1").Value2 = "基金名称"
$q3.Range("D1").Value2 = "基金规模"
$q3.Range("E1").Value2 = "股票总仓位"
$q3.Range("F1").Value2 = "仓位占比"
$q3.Range("G1").Value2 = "持有市值(亿元)"
$q3.Range("H1").Value2 = "仓位排名"

$q3.Range("A2").Value2 = 0
$q3.Range("B2").Value2 = "'160416"
$q3.Range("C2").Value2 = "华安标普全球石油指数（QDII-LOF）A"
$q3.Range("D2").Value2 = "'2.74"
$q3.Range("E2").Value2 = "'93.58"
$q3.Range("F2").Value2 = "'8.71"
$q3.Range("G2").Value2 = "'0.2387"
$q3.Range("H2").Value2 = 2

$q3.Range("A3").Value2 = 1
$q3.Range("B3").Value2 = "'014982"
$q3.Range("C3").Value2 = "华安标普全球石油指数（QDII-LOF）C"
$q3.Range("D3").Value2 = "'0.22"
$q3.Range("E3").Value2 = "'93.58"
$q3.Range("F3").Value2 = "'8.71"
$q3.Range("G3").Value2 = "'0.0192"
$q3.Range("H3").Value2 = 2

# ---------------------------------------------------------------------
# Step 3: add the 2022-Q3 summary row to the "总计" overview sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push existing rows 2..8 down to 3..9 (B/C/D values only, bottom-up so we
# never overwrite a row before it has been read). Column A is just a
# 0-based running index, so it is rewritten from scratch below instead of
# being shifted.
for ($r = 8; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    $total.Cells.Item($dst, 2).Value2 = $total.Cells.Item($src, 2).Value2
    $total.Cells.Item($dst, 3).Value2 = $total.Cells.Item($src, 3).Value2
    $total.Cells.Item($dst, 4).Value2 = $total.Cells.Item($src, 4).Value2
}

# Give the newly-exposed row 9 the same "A" column styling (bold, bordered)
# as the rest of the index column by copying the format down from row 8.
$total.Range("A8").Copy()
$total.Range("A9").PasteSpecial(-4122)

# Renumber the whole 0-based index column A2:A9.
for ($r = 2; $r -le 9; $r++) {
    $total.Cells.Item($r, 1).Value2 = $r - 2
}

# Write the new 2022-Q3 summary values into row 2.
$total.Range("B2").Value2 = "2022-Q3"
$total.Range("C2").Value2 = 2
$total.Range("D2").Value2 = 0.26

Write-Output "2022-Q3 sheet + summary row added"
